# Nexial "#system" reference sheet: add a new "localdb" command category.
#
# The #system sheet keeps, in column A ("target"), an alphabetically sorted
# list of command-type categories; each category's own column (starting at
# B) holds the list of function signatures that belong to it. A defined
# name (matching the category) marks the exact range used for the
# in-sheet autocomplete/validation list.
#
# "localdb" is a brand-new category that sorts alphabetically between
# "json" and "macro". Because "macro" already lived in column N, inserting
# a whole new column at N pushes "macro" (and every category after it)
# one column to the right, freeing up column N for "localdb"'s six
# function names. Column A also needs "localdb" spliced in at row 14
# (between "json" and "macro"), pushing every following category down by
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert a new column at N; this shifts the existing N:AC columns
#    (macro..xml) one column to the right, to O:AD, preserving their data.
$ws.Columns("N").Insert(-4161)

# 2) Populate the freshly emptied column N with the "localdb" category:
#    header in row 1, six function signatures in rows 2-7.
$ws.Range("N1").Value2 = "localdb"
$ws.Range("N2").Value2 = "cloneTable(var,source,target)"
$ws.Range("N3").Value2 = "dropTables(var,tables)"
$ws.Range("N4").Value2 = "exportCSV(sql,output)"
$ws.Range("N5").Value2 = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value2 = "purge(var)"
$ws.Range("N7").Value2 = "runSQLs(var,sqls)"

# 3) Splice "localdb" into the sorted category list in column A (row 14),
#    pushing "macro" and everything after it down by one row (A14:A29
#    become A15:A30). Column-A-only, so do it with direct value writes
#    rather than a row insert (which would shift every column). Note:
#    use Value2 (not Value) to read/write single-cell scalars here.
$ws.Range("A30").Value2 = $ws.Range("A29").Value2
$ws.Range("A29").Value2 = $ws.Range("A28").Value2
$ws.Range("A28").Value2 = $ws.Range("A27").Value2
$ws.Range("A27").Value2 = $ws.Range("A26").Value2
$ws.Range("A26").Value2 = $ws.Range("A25").Value2
$ws.Range("A25").Value2 = $ws.Range("A24").Value2
$ws.Range("A24").Value2 = $ws.Range("A23").Value2
$ws.Range("A23").Value2 = $ws.Range("A22").Value2
$ws.Range("A22").Value2 = $ws.Range("A21").Value2
$ws.Range("A21").Value2 = $ws.Range("A20").Value2
$ws.Range("A20").Value2 = $ws.Range("A19").Value2
$ws.Range("A19").Value2 = $ws.Range("A18").Value2
$ws.Range("A18").Value2 = $ws.Range("A17").Value2
$ws.Range("A17").Value2 = $ws.Range("A16").Value2
$ws.Range("A16").Value2 = $ws.Range("A15").Value2
$ws.Range("A15").Value2 = $ws.Range("A14").Value2
$ws.Range("A14").Value2 = "localdb"

# 4) Repoint the defined names whose target ranges moved because of the
#    column insert (one column to the right, N:AC -> O:AD) and the
#    "target" category list growing by one row.
$wb.Names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"

# 5) Add the new "localdb" defined name used for its autocomplete list.
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
